# 15.1.2.1.xlsx -- extend the indicator table with two more years (2022, 2023)
# and tidy up row heights / the stale cell selection left in the sheet view.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row heights -----------------------------------------------------
# Row 2 gains an explicit 15pt row height, row 3 grows from 12.75 to 15pt,
# and row 4 shrinks slightly from 18.75 to 17.25pt.
$ws.Rows("2:2").RowHeight = 15
$ws.Rows("3:3").RowHeight = 15
$ws.Rows("4:4").RowHeight = 17.25

# --- New data columns: N (2022) and O (2023) --------------------------
# Put the values in first, then copy the neighbouring (2021) column's
# formatting onto the new cells so they pick up the same styles as the
# rest of the header/data rows.
$ws.Range("N4").Value = 2022
$ws.Range("O4").Value = 2023
$ws.Range("N5").Value = 6.53
$ws.Range("O5").Value = 6.53

$ws.Range("M4").Copy()
$ws.Range("N4:O4").PasteSpecial(-4122)

$ws.Range("M5").Copy()
$ws.Range("N5:O5").PasteSpecial(-4122)

$excel.CutCopyMode = $false

# --- Clear the leftover cell selection ---------------------------------
# The workbook had been saved with P6 selected; reset the active cell
# back to the top of the sheet.
$ws.Range("A1").Select()
